$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal text, preserving its original style/format
function Set-TextValue($cellAddr, $value) {
    $cell = $ws.Range($cellAddr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = '29.585.06'
$ws.Range("E2").Value = '  +2.33%  '
$ws.Range("D3").Value = '1.858.03'
$ws.Range("E3").Value = '  +1.48%  '
Set-TextValue "D4" '0.9995'
$ws.Range("E4").Value = '  +0.00%  '
Set-TextValue "D5" '244.64'
$ws.Range("E5").Value = '  +0.17%  '
Set-TextValue "D6" '0.6937'
$ws.Range("E6").Value = '  +0.81%  '
Set-TextValue "D7" '1.000'
$ws.Range("E7").Value = '  +0.01%  '
Set-TextValue "D8" '0.07695'
Set-TextValue "D9" '0.3058'
$ws.Range("E9").Value = '  +0.09%  '
Set-TextValue "D10" '23.71'
$ws.Range("E10").Value = '  +0.73%  '
Set-TextValue "D11" '0.07772'
$ws.Range("E11").Value = '  -0.38%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D12" '5.147'
$ws.Range("E12").Value = '  +1.37%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.856.86'
$ws.Range("E13").Value = '  +1.37%  '
Set-TextValue "D14" '91.55'
$ws.Range("E14").Value = '  +1.14%  '
Set-TextValue "D15" '0.6917'
$ws.Range("E15").Value = '  +1.99%  '
Set-TextValue "D16" '6.563'
$ws.Range("E16").Value = '  +1.70%  '
$ws.Range("D17").Value = '29.565.68'
$ws.Range("E17").Value = '  +2.14%  '
Set-TextValue "D18" '0.000008277'
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("D19").Value = '2.103.46'
$ws.Range("E19").Value = '  +0.71%  '
Set-TextValue "D20" '240.17'
$ws.Range("E20").Value = '  -1.12%  '
$ws.Range("E21").Value = '  +0.63%  '
Set-TextValue "D22" '1.0000'
$ws.Range("E22").Value = '  +0.00%  '
Set-TextValue "D23" '7.594'
Set-TextValue "D24" '0.9998'
$ws.Range("E24").Value = '  +0.04%  '
Set-TextValue "D25" '0.1498'
$ws.Range("E25").Value = '  +1.79%  '
Set-TextValue "D26" '8.917'
$ws.Range("E26").Value = '  +1.44%  '
Set-TextValue "D27" '159.68'
$ws.Range("E27").Value = '  -1.07%  '
Set-TextValue "D28" '18.28'
$ws.Range("E28").Value = '  +0.56%  '
Set-TextValue "D29" '1.532'
$ws.Range("E29").Value = '  -1.06%  '
Set-TextValue "D30" '4.249'
$ws.Range("E30").Value = '  +0.83%  '
Set-TextValue "D31" '4.181'
$ws.Range("E31").Value = '  +1.30%  '
Set-TextValue "D32" '1.200'
$ws.Range("E32").Value = '  +1.84%  '
$ws.Range("E33").Value = '  -0.63%  '
Set-TextValue "D34" '0.7720'
$ws.Range("E34").Value = '  +2.10%  '
Set-TextValue "D35" '1.894'
$ws.Range("E35").Value = '  +3.26%  '
Set-TextValue "D36" '1.152'
$ws.Range("E36").Value = '  +0.67%  '
Set-TextValue "D37" '2.684'
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("D38").Value = '1.334.57'
$ws.Range("E38").Value = '  +8.46%  '
$ws.Range("E39").Value = '  +1.36%  '
$ws.Range("E40").Value = '  +1.13%  '
Set-TextValue "D41" '0.9692'
$ws.Range("E41").Value = '  +4.92%  '
Set-TextValue "D42" '106.32'
$ws.Range("E42").Value = '  -2.05%  '
Set-TextValue "D43" '5.781'
$ws.Range("E43").Value = '  +0.89%  '
Set-TextValue "D44" '1.000'
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D45" '9.766'
$ws.Range("E45").Value = '  +2.79%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '2.000.55'
$ws.Range("E46").Value = '  +1.01%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue "D47" '0.00000000123'
$ws.Range("E47").Value = '  +1.78%  '
Set-TextValue "D48" '0.5216'
$ws.Range("E48").Value = '  +0.94%  '
Set-TextValue "D49" '1.773'
Set-TextValue "D50" '63.50'
$ws.Range("E50").Value = '  -1.31%  '
$ws.Range("E51").Value = '  +0.74%  '
